$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of column J ("k" parameter) across rows 2-11, bold
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$jfont = $ws.Range("J12").Font
$jfont.Bold = $true
$jfont.Size = 11

# Rows 14-17: summary labels (column A) + aggregate formulas (column B)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Formatting for the summary values: bold, larger font, vertically centered
$rng = $ws.Range("B14:B17")
$rfont = $rng.Font
$rfont.Bold = $true
$rfont.Size = 12
$rng.VerticalAlignment = -4108

# Taller rows for the summary block
$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

# Selection matching the saved file: A14 active, A14:B17 selected
[void]$ws.Range("A14:B17").Select()

# Page setup: paper size 9 (A4), portrait orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
